$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 766.4286
$ws.Range("I107").Value = 571.5294
$ws.Range("K107").Value = 571.5294
$ws.Range("M107").Value = 1348.4706

$ws.Range("H127").Value = 7200
$ws.Range("I127").Value = 1000
$ws.Range("J127").Value = 11333.333
$ws.Range("K127").Value = 3000
$ws.Range("L127").Value = 33999.999
$ws.Range("M127").Value = 1960
$ws.Range("N127").Value = -43919.999

$ws.Range("H130").Value = 45666.668
$ws.Range("J130").Value = 45666.668
$ws.Range("L130").Value = 45666.668
$ws.Range("N130").Value = -55706.668

$ws.Range("H132").Value = 3748.353
$ws.Range("J132").Value = 14000
$ws.Range("L132").Value = 42000
$ws.Range("N132").Value = -47060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3961.2942
$ws.Range("I122").Value = 2766
$ws.Range("K122").Value = 8298
$ws.Range("M122").Value = -5848

$ws.Range("H132").Value = 8191.6895
$ws.Range("I132").Value = 4474.933
$ws.Range("K132").Value = 13424.799
$ws.Range("M132").Value = -10894.799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 241.33333
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 700
$ws.Range("K23").Value = 12
$ws.Range("L23").Value = 700
$ws.Range("M23").Value = 271
$ws.Range("N23").Value = -1266

$ws.Range("H63").Value = 8666.666999999999
$ws.Range("J63").Value = 8666.666999999999
$ws.Range("L63").Value = 8666.666999999999
$ws.Range("N63").Value = -10038.667

$ws.Range("H66").Value = 8666.666999999999
$ws.Range("J66").Value = 8666.666999999999
$ws.Range("L66").Value = 26000.001
$ws.Range("N66").Value = -32864.001

$ws.Range("H75").Value = 6999.5
$ws.Range("I75").Value = 6999.5
$ws.Range("K75").Value = 6999.5
$ws.Range("M75").Value = -6063.5

$ws.Range("H78").Value = 6999.5
$ws.Range("I78").Value = 6999.5
$ws.Range("K78").Value = 20998.5
$ws.Range("M78").Value = -16318.5

$ws.Range("H94").Value = 701.4211
$ws.Range("I94").Value = 664.3125
$ws.Range("K94").Value = 664.3125
$ws.Range("M94").Value = -213.3125

$ws.Range("H99").Value = 12544.182
$ws.Range("I99").Value = 21957.6
$ws.Range("K99").Value = 21957.6
$ws.Range("M99").Value = -20459.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 192281.55
$ws.Range("J22").Value = 5500
$ws.Range("L22").Value = 5500
$ws.Range("N22").Value = -6200

$ws.Range("H134").Value = 1253874.8
$ws.Range("I134").Value = 1253874.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3761624.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3759089.4
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5541.909
$ws.Range("I55").Value = 5995.75
$ws.Range("J55").Value = 5282.5713
$ws.Range("K55").Value = 17987.25
$ws.Range("L55").Value = 15847.7139
$ws.Range("M55").Value = -17810.25
$ws.Range("N55").Value = -16201.7139

$ws.Range("H56").Value = 7173.478
$ws.Range("I56").Value = 7173.478
$ws.Range("K56").Value = 7173.478
$ws.Range("M56").Value = -6643.478

$ws.Range("H68").Value = 1891.9584
$ws.Range("J68").Value = 1845.4
$ws.Range("L68").Value = 5536.200000000001
$ws.Range("N68").Value = -7158.200000000001

$ws.Range("H71").Value = 1891.9584
$ws.Range("J71").Value = 1845.4
$ws.Range("L71").Value = 16608.6
$ws.Range("N71").Value = -24720.6

$ws.Range("H109").Value = 1459.5
$ws.Range("I109").Value = 951.4
$ws.Range("K109").Value = 2854.2
$ws.Range("M109").Value = -1814.2

$ws.Range("H113").Value = 1965.8
$ws.Range("I113").Value = 1544
$ws.Range("K113").Value = 4632
$ws.Range("M113").Value = -2462

$ws.Range("H132").Value = 2215.7
$ws.Range("I132").Value = 2266.8125
$ws.Range("J132").Value = 2011.25
$ws.Range("K132").Value = 20401.3125
$ws.Range("L132").Value = 18101.25
$ws.Range("M132").Value = -17871.3125
$ws.Range("N132").Value = -23161.25

$ws.Range("H141").Value = 102097.06
$ws.Range("I141").Value = 115000.445
$ws.Range("K141").Value = 345001.335
$ws.Range("M141").Value = -339821.335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 36799.555
$ws.Range("I44").Value = 36399.5
$ws.Range("K44").Value = 36399.5
$ws.Range("M44").Value = -35803.5

$ws.Range("H53").Value = 5000
$ws.Range("I53").Value = 5000
$ws.Range("K53").Value = 5000
$ws.Range("M53").Value = -4369

$ws.Range("H70").Value = 8000
$ws.Range("I70").Value = 8000
$ws.Range("K70").Value = 8000
$ws.Range("M70").Value = -7730

$ws.Range("H73").Value = 8000
$ws.Range("I73").Value = 8000
$ws.Range("K73").Value = 8000
$ws.Range("M73").Value = -7064

$ws.Range("H97").Value = 1490.8334
$ws.Range("I97").Value = 1554.15
$ws.Range("J97").Value = 1174.25
$ws.Range("K97").Value = 1554.15
$ws.Range("L97").Value = 1174.25
$ws.Range("M97").Value = -1058.15
$ws.Range("N97").Value = -2166.25

$ws.Range("H130").Value = 159499.75
$ws.Range("J130").Value = 159499.75
$ws.Range("L130").Value = 159499.75
$ws.Range("N130").Value = -169539.75

$ws.Range("H141").Value = 108333
$ws.Range("J141").Value = 108333
$ws.Range("L141").Value = 108333
$ws.Range("N141").Value = -118693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 171653.67
$ws.Range("I7").Value = 5952.5
$ws.Range("K7").Value = 5952.5
$ws.Range("M7").Value = -5840.5

$ws.Range("H16").Value = 1417.5
$ws.Range("I16").Value = 1414.125
$ws.Range("K16").Value = 1414.125
$ws.Range("M16").Value = -1244.125

$ws.Range("H61").Value = 2392.625
$ws.Range("I61").Value = 2465.2856
$ws.Range("K61").Value = 2465.2856
$ws.Range("M61").Value = -2263.2856

$ws.Range("H93").Value = 76924424
$ws.Range("I93").Value = 83334380
$ws.Range("K93").Value = 83334380
$ws.Range("M93").Value = -83333132

$ws.Range("H113").Value = 2392.625
$ws.Range("I113").Value = 2465.2856
$ws.Range("K113").Value = 2465.2856
$ws.Range("M113").Value = -295.2856000000002

$ws.Range("H122").Value = 6753.231
$ws.Range("I122").Value = 5905
$ws.Range("K122").Value = 17715
$ws.Range("M122").Value = -15265

$ws.Range("H126").Value = 171653.67
$ws.Range("I126").Value = 5952.5
$ws.Range("K126").Value = 17857.5
$ws.Range("M126").Value = -15387.5

$ws.Range("H133").Value = 100499.5
$ws.Range("J133").Value = 100499.5
$ws.Range("L133").Value = 100499.5
$ws.Range("N133").Value = -105559.5

$ws.Range("H137").Value = 70199
$ws.Range("J137").Value = 70199
$ws.Range("L137").Value = 70199
$ws.Range("N137").Value = -80399

$ws.Range("H141").Value = 262500
$ws.Range("J141").Value = 262500
$ws.Range("L141").Value = 262500
$ws.Range("N141").Value = -272860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H46").Value = 87136.71000000001
$ws.Range("J46").Value = 87136.71000000001
$ws.Range("L46").Value = 87136.71000000001
$ws.Range("N46").Value = -87598.71000000001

$ws.Range("H54").Value = 23813
$ws.Range("I54").Value = 19035
$ws.Range("J54").Value = 26998.334
$ws.Range("K54").Value = 19035
$ws.Range("L54").Value = 26998.334
$ws.Range("M54").Value = -18515
$ws.Range("N54").Value = -28038.334

$ws.Range("H81").Value = 37723.668
$ws.Range("I81").Value = 32812.715
$ws.Range("K81").Value = 65625.42999999999
$ws.Range("M81").Value = -64564.42999999999

$ws.Range("H84").Value = 37723.668
$ws.Range("I84").Value = 32812.715
$ws.Range("K84").Value = 328127.15
$ws.Range("M84").Value = -322823.15

$ws.Range("H112").Value = 82500
$ws.Range("J112").Value = 82500
$ws.Range("L112").Value = 82500
$ws.Range("N112").Value = -85454

$ws.Range("H122").Value = 8709.875
$ws.Range("I122").Value = 6382.857
$ws.Range("K122").Value = 19148.571
$ws.Range("M122").Value = -16698.571

$ws.Range("H134").Value = 87136.71000000001
$ws.Range("J134").Value = 87136.71000000001
$ws.Range("L134").Value = 261410.13
$ws.Range("N134").Value = -266480.13
